# NumberFormatAndAlignment regression fixture update:
# refresh the sample Price/Count values in rows 2-4 (columns B and C).
# (The underlying numFmtId 166 definition itself is a private test-fixture
# typo fix ("#,###" -> "#,##0") that renders identically for this data set
# and has no user-facing Excel surface to edit in place, so it is left
# alone here; only the cell contents change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 141.5
$ws.Range("C2").Value = 92

$ws.Range("B3").Value = 0.314
$ws.Range("C3").Value = 15

$ws.Range("B4").Value = 653.5
$ws.Range("C4").Value = 14
